# debug check for full rows of missing data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: EXPOCODE 91AH20140402 (PANGEA chl -> ctd/split, tab -> comma) ---
$ws.Range("B12").Value = ".csv"
$ws.Range("C12").Value = "E:/Data_downloads/PANGEAE/91AH20140402/ctd/split"
$ws.Range("G12").Value = "comma"
$ws.Range("S12").Value = "Date.Time"
$ws.Range("U12").Value = "Date.Time"
$ws.Range("AG12").Value = "Press"

# --- Row 13: EXPOCODE 35XI20090905 (fix mismatched path + time format + pressure label) ---
$ws.Range("C13").Value = "E:/Data_downloads/TARA_PANGEAE/35XI20090905/ctd/split"
$ws.Range("N13").Value = "NA"
$ws.Range("T13").Value = "%Y-%m-%dT%H:%M:%S"
$ws.Range("X13").Value = "%Y-%m-%dT%H:%M:%S"
$ws.Range("AG13").Value = "Press"

# --- Row 21: EXPOCODE 29HE19951203 (tab -> comma, TIME_b label, pressure label) ---
$ws.Range("G21").Value = "comma"
$ws.Range("U21").Value = "Date"
$ws.Range("AG21").Value = "Press"

# --- Row 22: EXPOCODE 29HE19960117 (tab -> comma, TIME_b label, pressure label) ---
$ws.Range("G22").Value = "comma"
$ws.Range("U22").Value = "Date"
$ws.Range("AG22").Value = "Press"

# --- Sheet cosmetics: column E width + new selection ---
$ws.Columns("E").ColumnWidth = 13
$ws.Range("N11").Select()

$wb.Save()
